# Fix a bug in targetScene: the data rows (A:F) were associated with the
# wrong "scene" id. The fix re-orders rows 2-25 (the data rows, excluding
# the totals row 26) so that each row of values A:F ends up on the correct
# row, without altering any of the values themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current contents of the data range (A2:F25) *before* writing
# anything back, since several rows swap places with each other.
$source = $ws.Range("A2:F25").Value2

# Mapping of old row number -> new row number (1-based, relative to the
# whole sheet). Rows 10, 16, 17, 22-25 (and the header/totals rows) keep
# their original position.
$rowMap = @{
    2 = 3
    3 = 8
    4 = 13
    5 = 15
    6 = 12
    7 = 5
    8 = 2
    9 = 14
    10 = 10
    11 = 4
    12 = 6
    13 = 9
    14 = 11
    15 = 7
    16 = 16
    17 = 17
    18 = 19
    19 = 21
    20 = 18
    21 = 20
    22 = 22
    23 = 23
    24 = 24
    25 = 25
}

foreach ($oldRow in $rowMap.Keys) {
    $newRow = $rowMap[$oldRow]
    $srcIndex = $oldRow - 1   # offset into the A2:F25 snapshot (row 2 -> index 1)
    for ($col = 1; $col -le 6; $col++) {
        $value = $source[$srcIndex, $col]
        $ws.Cells.Item($newRow, $col).Value2 = $value
    }
}
